$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet: Property1 -> DataNode (unifying DataNode/DataTable/Entity concepts)
$ws.Name = "DataNode"

# Row height tweaks (header row + wrapped description row), matching the
# Windows-Excel re-save of the workbook.
$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 40.5

# Restore the last-used selection from the author's editing session.
$ws.Range("C36").Select()
